$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1292.3334
$ws.Range("J2").Value = 868.25
$ws.Range("L2").Value = 868.25
$ws.Range("N2").Value = -1094.25

$ws.Range("H4").Value = 429.14285
$ws.Range("I4").Value = 412.8
$ws.Range("K4").Value = 412.8
$ws.Range("M4").Value = -298.8

$ws.Range("H33").Value = 290.26666
$ws.Range("I33").Value = 201.9
$ws.Range("K33").Value = 201.9
$ws.Range("M33").Value = 27.09999999999999

$ws.Range("H43").Value = 5750
$ws.Range("J43").Value = 5857
$ws.Range("L43").Value = 5857
$ws.Range("N43").Value = -5995

$ws.Range("H88").Value = 7714.364
$ws.Range("J88").Value = 8385.799999999999
$ws.Range("L88").Value = 8385.799999999999
$ws.Range("N88").Value = -9197.799999999999

$ws.Range("H91").Value = 7714.364
$ws.Range("J91").Value = 8385.799999999999
$ws.Range("L91").Value = 8385.799999999999
$ws.Range("N91").Value = -11193.8

$ws.Range("H96").Value = 5102923
$ws.Range("I96").Value = 8928902
$ws.Range("J96").Value = 1618.3334
$ws.Range("K96").Value = 26786706
$ws.Range("L96").Value = 4855.0002
$ws.Range("M96").Value = -26785333
$ws.Range("N96").Value = -7601.0002

$ws.Range("H116").Value = 45272.668
$ws.Range("I116").Value = 45272.668
$ws.Range("K116").Value = 45272.668
$ws.Range("M116").Value = -41830.668

$ws.Range("H131").Value = 2848
$ws.Range("I131").Value = 2304.7144
$ws.Range("J131").Value = 4749.5
$ws.Range("K131").Value = 6914.1432
$ws.Range("L131").Value = 14248.5
$ws.Range("M131").Value = -1874.1432
$ws.Range("N131").Value = -24328.5

$ws.Range("H137").Value = 25000
$ws.Range("I137").Value = 37089.445
$ws.Range("K137").Value = 111268.335
$ws.Range("M137").Value = -108718.335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 792.75
$ws.Range("I2").Value = 673.8333
$ws.Range("K2").Value = 673.8333
$ws.Range("M2").Value = -560.8333

$ws.Range("H45").Value = 2943.5557
$ws.Range("I45").Value = 1718.2727
$ws.Range("J45").Value = 4869
$ws.Range("K45").Value = 1718.2727
$ws.Range("L45").Value = 4869
$ws.Range("M45").Value = -1341.2727
$ws.Range("N45").Value = -5623

$ws.Range("H97").Value = 2187.375
$ws.Range("I97").Value = 1938.8462
$ws.Range("K97").Value = 1938.8462
$ws.Range("M97").Value = -1442.8462

$ws.Range("H116").Value = 792.75
$ws.Range("I116").Value = 673.8333
$ws.Range("K116").Value = 673.8333
$ws.Range("M116").Value = 1620.1667

$ws.Range("H122").Value = 2948.7778
$ws.Range("I122").Value = 2769
$ws.Range("J122").Value = 3173.5
$ws.Range("K122").Value = 8307
$ws.Range("L122").Value = 9520.5
$ws.Range("M122").Value = -5857
$ws.Range("N122").Value = -14420.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 792.75
$ws.Range("I3").Value = 673.8333
$ws.Range("K3").Value = 673.8333
$ws.Range("M3").Value = -559.8333

$ws.Range("H86").Value = 3562.2856
$ws.Range("I86").Value = 2367.3
$ws.Range("J86").Value = 6549.75
$ws.Range("K86").Value = 2367.3
$ws.Range("L86").Value = 6549.75
$ws.Range("M86").Value = -1244.3
$ws.Range("N86").Value = -8795.75

$ws.Range("H89").Value = 3562.2856
$ws.Range("I89").Value = 2367.3
$ws.Range("J89").Value = 6549.75
$ws.Range("K89").Value = 11836.5
$ws.Range("L89").Value = 32748.75
$ws.Range("M89").Value = -6220.5
$ws.Range("N89").Value = -43980.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12504238
$ws.Range("I31").Value = 33333998
$ws.Range("J31").Value = 6382.4
$ws.Range("K31").Value = 33333998
$ws.Range("L31").Value = 6382.4
$ws.Range("M31").Value = -33333703
$ws.Range("N31").Value = -6972.4

$ws.Range("H34").Value = 12504238
$ws.Range("I34").Value = 33333998
$ws.Range("J34").Value = 6382.4
$ws.Range("K34").Value = 33333998
$ws.Range("L34").Value = 6382.4
$ws.Range("M34").Value = -33333796
$ws.Range("N34").Value = -6786.4

$ws.Range("H47").Value = 10000
$ws.Range("I47").Value = 10000
$ws.Range("K47").Value = 10000
$ws.Range("M47").Value = -9434

$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws.Range("H134").Value = 4372.75
$ws.Range("I134").Value = 4196.4
$ws.Range("K134").Value = 12589.2
$ws.Range("M134").Value = -10054.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 41972.555
$ws.Range("J37").Value = 41972.555
$ws.Range("L37").Value = 125917.665
$ws.Range("N37").Value = -126141.665

$ws.Range("H68").Value = 4107.5
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 4211.8
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 12635.4
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -14257.4

$ws.Range("H71").Value = 4107.5
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 4211.8
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 37906.2
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -46018.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 601.6842
$ws.Range("I2").Value = 795
$ws.Range("J2").Value = 270.2857
$ws.Range("K2").Value = 795
$ws.Range("L2").Value = 270.2857
$ws.Range("M2").Value = -682
$ws.Range("N2").Value = -496.2857

$ws.Range("H32").Value = 26250
$ws.Range("J32").Value = 26250
$ws.Range("L32").Value = 26250
$ws.Range("N32").Value = -26842

$ws.Range("H42").Value = 37899
$ws.Range("J42").Value = 37899
$ws.Range("L42").Value = 37899
$ws.Range("N42").Value = -38869

$ws.Range("H80").Value = 14932.833
$ws.Range("I80").Value = 7765.6665
$ws.Range("J80").Value = 22100
$ws.Range("K80").Value = 7765.6665
$ws.Range("L80").Value = 22100
$ws.Range("M80").Value = -6767.6665
$ws.Range("N80").Value = -24096

$ws.Range("H83").Value = 14932.833
$ws.Range("I83").Value = 7765.6665
$ws.Range("J83").Value = 22100
$ws.Range("K83").Value = 38828.3325
$ws.Range("L83").Value = 110500
$ws.Range("M83").Value = -33836.3325
$ws.Range("N83").Value = -120484

$ws.Range("H113").Value = 1306.6
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

$ws.Range("H115").Value = 37899
$ws.Range("J115").Value = 37899
$ws.Range("L115").Value = 37899
$ws.Range("N115").Value = -40249

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3913.5715
$ws.Range("I7").Value = 3913.5715
$ws.Range("K7").Value = 3913.5715
$ws.Range("M7").Value = -3801.5715

$ws.Range("H46").Value = 4714.636
$ws.Range("I46").Value = 2266.6667
$ws.Range("K46").Value = 2266.6667
$ws.Range("M46").Value = -2078.6667

$ws.Range("H100").Value = 3316.5
$ws.Range("I100").Value = 2866.3333
$ws.Range("J100").Value = 3766.6667
$ws.Range("K100").Value = 2866.3333
$ws.Range("L100").Value = 3766.6667
$ws.Range("M100").Value = -2325.3333
$ws.Range("N100").Value = -4848.6667

$ws.Range("H102").Value = 70000
$ws.Range("J102").Value = 70000
$ws.Range("L102").Value = 70000
$ws.Range("N102").Value = -76490

$ws.Range("H126").Value = 3913.5715
$ws.Range("I126").Value = 3913.5715
$ws.Range("K126").Value = 11740.7145
$ws.Range("M126").Value = -9270.7145

$ws.Range("H132").Value = 2165
$ws.Range("I132").Value = 1679.7273
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 5039.1819
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -2509.1819
$ws.Range("N132").Value = -15558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -43134

$ws.Range("H122").Value = 43232.2
$ws.Range("I122").Value = 49671.633
$ws.Range("K122").Value = 149014.899
$ws.Range("M122").Value = -146564.899

$ws.Range("H132").Value = 871.4545000000001
$ws.Range("I132").Value = 730.4286
$ws.Range("K132").Value = 2191.2858
$ws.Range("M132").Value = 338.7142000000003
